# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 159 in the Mandarina sheet,
# shifting the existing rows 159:227 down to 160:228.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 159 (pushes old row 159.. down by one).
$ws.Rows(159).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(159, 1).Value  = 11
$ws.Cells.Item(159, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(159, 3).Value  = "Bíobío"
$ws.Cells.Item(159, 4).Value  = 45119
$ws.Cells.Item(159, 5).Value  = 8
$ws.Cells.Item(159, 6).Value  = "Fruta"
$ws.Cells.Item(159, 7).Value  = 100102
$ws.Cells.Item(159, 8).Value  = "Cítricos"
$ws.Cells.Item(159, 9).Value  = 100102004
$ws.Cells.Item(159, 10).Value = "Mandarina"
$ws.Cells.Item(159, 11).Value = "Clementina"
$ws.Cells.Item(159, 12).Value = "Primera"
$ws.Cells.Item(159, 13).Value = 200
$ws.Cells.Item(159, 14).Value = 9000
$ws.Cells.Item(159, 15).Value = 10000
$ws.Cells.Item(159, 16).Value = 9500
$ws.Cells.Item(159, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(159, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(159, 19).Value = 950
$ws.Cells.Item(159, 20).Value = 10
